$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 49

# Column A holds a literal date-formatted text string (e.g. "08/08/2025" in the
# row directly above), not a real date value. Assigning a date-looking string
# straight to .Value would get auto-converted into a date serial number, so
# format the cell as Text first, put the literal string in, then drop the
# cell back to the default "Normal" style (the other text-date rows carry no
# explicit style either) while keeping its contents as plain text.
$ws.Cells.Item($row, 1).NumberFormat = "@"
$ws.Cells.Item($row, 1).Value = "08/10/2025"
$ws.Cells.Item($row, 1).Style = "Normal"

$ws.Cells.Item($row, 2).Value = 107.096000000005
$ws.Cells.Item($row, 3).Value = 0.0933741689698918
$ws.Cells.Item($row, 4).Value = 10
